$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number-format on the price (D) cells being updated so that
# numeric-looking strings (e.g. "16.64", "1.000", "11.00") are stored as
# literal text instead of being coerced to numbers (which would drop
# formatting such as trailing zeros). NumberFormat is set cell-by-cell
# (a union Range(...) only honors the first area for this property).
$priceCells = @("D2", "D3", "D5", "D7", "D9", "D10", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D21", "D22", "D23", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "27.019.11"
$ws.Range("E2").Value = "  -2.92%  "
$ws.Range("D3").Value = "1.728.35"
$ws.Range("E3").Value = "  -1.47%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "310.21"
$ws.Range("E5").Value = "  -5.14%  "
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("D7").Value = "0.4842"
$ws.Range("E7").Value = "  +3.95%  "
$ws.Range("E8").Value = "  +0.40%  "
$ws.Range("D9").Value = "43.37"
$ws.Range("E9").Value = "  +3.38%  "
$ws.Range("D10").Value = "0.07260"
$ws.Range("E10").Value = "  -1.19%  "
$ws.Range("E11").Value = "  -2.60%  "
$ws.Range("E12").Value = "  +0.04%  "
$ws.Range("D13").Value = "19.95"
$ws.Range("E13").Value = "  -2.53%  "
$ws.Range("D14").Value = "5.883"
$ws.Range("E14").Value = "  -1.46%  "
$ws.Range("D15").Value = "1.726.61"
$ws.Range("E15").Value = "  -1.47%  "
$ws.Range("D16").Value = "6.866"
$ws.Range("E16").Value = "  -3.70%  "
$ws.Range("D17").Value = "87.18"
$ws.Range("E17").Value = "  -5.21%  "
$ws.Range("D18").Value = "0.00001036"
$ws.Range("E18").Value = "  -1.47%  "
$ws.Range("D19").Value = "0.06407"
$ws.Range("E19").Value = "  +0.03%  "
$ws.Range("E20").Value = "  +0.07%  "
$ws.Range("D21").Value = "16.64"
$ws.Range("E21").Value = "  -0.81%  "
$ws.Range("D22").Value = "5.712"
$ws.Range("E22").Value = "  -0.60%  "
$ws.Range("D23").Value = "27.089.13"
$ws.Range("E23").Value = "  -2.79%  "
$ws.Range("E24").Value = "  -1.86%  "
$ws.Range("D25").Value = "2.077"
$ws.Range("E25").Value = "  -3.71%  "
$ws.Range("D26").Value = "154.49"
$ws.Range("E26").Value = "  -4.50%  "
$ws.Range("D27").Value = "19.97"
$ws.Range("E27").Value = "  -0.19%  "
$ws.Range("D28").Value = "1.923.03"
$ws.Range("E28").Value = "  -1.62%  "
$ws.Range("D29").Value = "2.075"
$ws.Range("E29").Value = "  -2.88%  "
$ws.Range("D30").Value = "120.91"
$ws.Range("E30").Value = "  -1.34%  "
$ws.Range("D31").Value = "1.052"
$ws.Range("E31").Value = "  -0.92%  "
$ws.Range("D32").Value = "0.09354"
$ws.Range("E32").Value = "  +0.76%  "
$ws.Range("D33").Value = "3.654"
$ws.Range("E33").Value = "  -0.02%  "
$ws.Range("D34").Value = "5.376"
$ws.Range("E34").Value = "  -2.67%  "
$ws.Range("D35").Value = "0.05963"
$ws.Range("E35").Value = "  -1.72%  "
$ws.Range("D36").Value = "0.02186"
$ws.Range("E36").Value = "  -3.49%  "
$ws.Range("B37").Value = "Aptos"
$ws.Range("C37").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D37").Value = "11.00"
$ws.Range("E37").Value = "  -5.17%  "
$ws.Range("B38").Value = "WEMIXTOKEN"
$ws.Range("C38").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D38").Value = "1.428"
$ws.Range("E38").Value = "  +6.34%  "
$ws.Range("B39").Value = "Algorand"
$ws.Range("C39").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D39").Value = "0.1995"
$ws.Range("E39").Value = "  -3.07%  "
$ws.Range("B40").Value = "InternetComputer(DFINITY)"
$ws.Range("C40").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D40").Value = "4.765"
$ws.Range("E40").Value = "  -2.51%  "
$ws.Range("D41").Value = "1.001"
$ws.Range("E41").Value = "  +0.66%  "
$ws.Range("D42").Value = "0.5986"
$ws.Range("E42").Value = "  -2.72%  "
$ws.Range("D44").Value = "7.508"
$ws.Range("E44").Value = "  -3.01%  "
$ws.Range("D45").Value = "12.76"
$ws.Range("E45").Value = "  -2.03%  "
$ws.Range("D46").Value = "3.579"
$ws.Range("E46").Value = "  -4.02%  "
$ws.Range("D47").Value = "0.5630"
$ws.Range("E47").Value = "  -2.28%  "
$ws.Range("D48").Value = "118.78"
$ws.Range("E48").Value = "  -3.05%  "
$ws.Range("D49").Value = "1.849"
$ws.Range("E49").Value = "  -3.77%  "
$ws.Range("D50").Value = "1.105"
$ws.Range("E50").Value = "  -1.15%  "
$ws.Range("D51").Value = "0.06644"
$ws.Range("E51").Value = "  -2.17%  "
